# example.xlsx now specifies paths through env
# Filepaths are specified through the environment variable GITHUB_WORKSPACE,
# allowing tests to be run locally or through GitHub Actions.

$wb = $excel.ActiveWorkbook

# "Samples" sheet: update the two path cells to use $GITHUB_WORKSPACE instead
# of the author's local / SharePoint-synced OneDrive path.
$samples = $wb.Sheets.Item("Samples")
$samples.Range("E2").Value = "`$GITHUB_WORKSPACE/test/inputs/pr_folder"
$samples.Range("E3").Value = "`$GITHUB_WORKSPACE/test/inputs/small.fcs"

# Move the active sheet/selection back to "Samples" (was left on
# "Transformations" with B3 selected; now Samples is active with E4 selected).
$samples.Activate()
$samples.Range("E4").Select()
